# Auto-generated Excel COM-interop script to apply the scheduled-runner value updates
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) of the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1216.6666
$ws.Range("J17").Value = 1159.45
$ws.Range("L17").Value = 3478.35
$ws.Range("N17").Value = -3814.35
$ws.Range("H40").Value = 8557.143
$ws.Range("I40").Value = 9900
$ws.Range("J40").Value = 8333.333000000001
$ws.Range("K40").Value = 9900
$ws.Range("L40").Value = 8333.333000000001
$ws.Range("M40").Value = -9725
$ws.Range("N40").Value = -8683.333000000001
$ws.Range("H62").Value = 1621.7778
$ws.Range("I62").Value = 1621.7778
$ws.Range("K62").Value = 1621.7778
$ws.Range("M62").Value = -997.7778000000001
$ws.Range("H65").Value = 1621.7778
$ws.Range("I65").Value = 1621.7778
$ws.Range("K65").Value = 8108.889
$ws.Range("M65").Value = -4988.889
$ws.Range("H111").Value = 3273.2856
$ws.Range("I111").Value = 2205
$ws.Range("J111").Value = 5944
$ws.Range("K111").Value = 6615
$ws.Range("L111").Value = 17832
$ws.Range("M111").Value = -3548
$ws.Range("N111").Value = -23966
$ws.Range("H118").Value = 1595.7142
$ws.Range("I118").Value = 1917.5
$ws.Range("J118").Value = 1166.6666
$ws.Range("K118").Value = 5752.5
$ws.Range("L118").Value = 3499.9998
$ws.Range("M118").Value = -4095.5
$ws.Range("N118").Value = -6813.9998
$ws.Range("H138").Value = 3366.395
$ws.Range("I138").Value = 1308.4722
$ws.Range("K138").Value = 3925.4166
$ws.Range("M138").Value = 1214.5834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15388801
$ws.Range("I32").Value = 15877301
$ws.Range("K32").Value = 15877301
$ws.Range("M32").Value = -15877014
$ws.Range("H61").Value = 3479.8572
$ws.Range("I61").Value = 3479.8572
$ws.Range("K61").Value = 3479.8572
$ws.Range("M61").Value = -3267.8572
$ws.Range("H74").Value = 2584.1052
$ws.Range("I74").Value = 2449.889
$ws.Range("K74").Value = 2449.889
$ws.Range("M74").Value = -1575.889
$ws.Range("H77").Value = 2584.1052
$ws.Range("I77").Value = 2449.889
$ws.Range("K77").Value = 12249.445
$ws.Range("M77").Value = -7881.445
$ws.Range("H88").Value = 2164.3333
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2496.5
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2496.5
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3308.5
$ws.Range("H91").Value = 2164.3333
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2496.5
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2496.5
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5304.5
$ws.Range("H102").Value = 1427
$ws.Range("I102").Value = 1427
$ws.Range("K102").Value = 1427
$ws.Range("M102").Value = 195
$ws.Range("H106").Value = 50690.668
$ws.Range("J106").Value = 50690.668
$ws.Range("L106").Value = 50690.668
$ws.Range("N106").Value = -53214.668
$ws.Range("H110").Value = 4556.222
$ws.Range("I110").Value = 4556.222
$ws.Range("K110").Value = 4556.222
$ws.Range("M110").Value = -2511.222
$ws.Range("H136").Value = 3479.8572
$ws.Range("I136").Value = 3479.8572
$ws.Range("K136").Value = 10439.5716
$ws.Range("M136").Value = -7889.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 2500
$ws.Range("K105").Value = 2500
$ws.Range("M105").Value = -753
$ws.Range("H134").Value = 1386.7333
$ws.Range("I134").Value = 1259.159
$ws.Range("K134").Value = 3777.477
$ws.Range("M134").Value = -1242.477

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8341.333000000001
$ws.Range("I16").Value = 8651.429
$ws.Range("K16").Value = 8651.429
$ws.Range("M16").Value = -8364.429
$ws.Range("H28").Value = 10000
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 1728.4
$ws.Range("I31").Value = 1625.2703
$ws.Range("K31").Value = 1625.2703
$ws.Range("M31").Value = -1330.2703
$ws.Range("H34").Value = 1728.4
$ws.Range("I34").Value = 1625.2703
$ws.Range("K34").Value = 1625.2703
$ws.Range("M34").Value = -1423.2703
$ws.Range("H43").Value = 38531.4
$ws.Range("J43").Value = 38531.4
$ws.Range("L43").Value = 38531.4
$ws.Range("N43").Value = -38899.4
$ws.Range("H62").Value = 5823
$ws.Range("I62").Value = 6506.3335
$ws.Range("K62").Value = 6506.3335
$ws.Range("M62").Value = -5882.3335
$ws.Range("H65").Value = 5823
$ws.Range("I65").Value = 6506.3335
$ws.Range("K65").Value = 32531.6675
$ws.Range("M65").Value = -29411.6675
$ws.Range("H101").Value = 38531.4
$ws.Range("J101").Value = 38531.4
$ws.Range("L101").Value = 38531.4
$ws.Range("N101").Value = -45021.4
$ws.Range("H107").Value = 13033.883
$ws.Range("I107").Value = 1303.5
$ws.Range("K107").Value = 1303.5
$ws.Range("M107").Value = 616.5
$ws.Range("H113").Value = 8341.333000000001
$ws.Range("I113").Value = 8651.429
$ws.Range("K113").Value = 8651.429
$ws.Range("M113").Value = -6481.429
$ws.Range("H134").Value = 2715.6296
$ws.Range("I134").Value = 2169.5386
$ws.Range("K134").Value = 6508.6158
$ws.Range("M134").Value = -3973.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 724.75
$ws.Range("I22").Value = 724.75
$ws.Range("K22").Value = 2174.25
$ws.Range("M22").Value = -2005.25
$ws.Range("H27").Value = 724.75
$ws.Range("I27").Value = 724.75
$ws.Range("K27").Value = 2174.25
$ws.Range("M27").Value = -2072.25
$ws.Range("H60").Value = 850
$ws.Range("I60").Value = 525
$ws.Range("J60").Value = 1175
$ws.Range("K60").Value = 1575
$ws.Range("L60").Value = 3525
$ws.Range("M60").Value = -1324
$ws.Range("N60").Value = -4027

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 34943.43
$ws.Range("J20").Value = 34943.43
$ws.Range("L20").Value = 34943.43
$ws.Range("N20").Value = -35433.43
$ws.Range("H107").Value = 931.75
$ws.Range("I107").Value = 575.6667
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 575.6667
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1344.3333
$ws.Range("N107").Value = -5840
$ws.Range("H122").Value = 7120.25
$ws.Range("J122").Value = 7201.4
$ws.Range("L122").Value = 21604.2
$ws.Range("N122").Value = -26504.2
$ws.Range("H126").Value = 5664.125
$ws.Range("I126").Value = 5151.5
$ws.Range("K126").Value = 15454.5
$ws.Range("M126").Value = -12984.5
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8713.294
$ws.Range("I93").Value = 8709.333000000001
$ws.Range("K93").Value = 8709.333000000001
$ws.Range("M93").Value = -7461.333000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H113").Value = 1431.5714
$ws.Range("I113").Value = 1239
$ws.Range("J113").Value = 1913
$ws.Range("K113").Value = 3717
$ws.Range("L113").Value = 5739
$ws.Range("M113").Value = -1547
$ws.Range("N113").Value = -10079
$ws.Range("H136").Value = 1423.9354
$ws.Range("I136").Value = 813.38464
$ws.Range("J136").Value = 4598.8
$ws.Range("K136").Value = 2440.15392
$ws.Range("L136").Value = 13796.4
$ws.Range("M136").Value = 109.8460800000003
$ws.Range("N136").Value = -18896.4

Write-Output "Updated 195 cells, cleared 4 cells across 8 sheets."